# Applies the cell-value updates (Betfair back/lay odds) described in the commit diff.
# Only numeric odds cells change; no rows/columns are inserted or removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 7).Value = 4.2  # G3: 4.7 -> 4.2
$ws.Cells.Item(3, 11).Value = 3.6  # K3: 3.65 -> 3.6
$ws.Cells.Item(3, 12).Value = 1.37  # L3: 1.01 -> 1.37
$ws.Cells.Item(3, 13).Value = 1.06  # M3: 1.07 -> 1.06
$ws.Cells.Item(3, 15).Value = 1.06  # O3: 1.31 -> 1.06
$ws.Cells.Item(3, 20).Value = 1.8  # T3: 1.78 -> 1.8
$ws.Cells.Item(3, 22).Value = 1.72  # V3: 1.71 -> 1.72
$ws.Cells.Item(3, 23).Value = 1.32  # W3: 1.27 -> 1.32
$ws.Cells.Item(3, 26).Value = 17.5  # Z3: 20 -> 17.5
$ws.Cells.Item(3, 30).Value = 14.5  # AD3: 16 -> 14.5
$ws.Cells.Item(3, 35).Value = 60  # AI3: 1000 -> 60
$ws.Cells.Item(3, 41).Value = 980  # AO3: 1000 -> 980
# Row 4
$ws.Cells.Item(4, 7).Value = 2.18  # G4: 2.24 -> 2.18
$ws.Cells.Item(4, 8).Value = 3.65  # H4: 3.5 -> 3.65
$ws.Cells.Item(4, 9).Value = 3.85  # I4: 3.8 -> 3.85
$ws.Cells.Item(4, 10).Value = 3.7  # J4: 3.6 -> 3.7
$ws.Cells.Item(4, 11).Value = 3.75  # K4: 3.7 -> 3.75
$ws.Cells.Item(4, 13).Value = 1.05  # M4: 1.01 -> 1.05
$ws.Cells.Item(4, 14).Value = 3.2  # N4: 2.14 -> 3.2
$ws.Cells.Item(4, 17).Value = 1.82  # Q4: 1.77 -> 1.82
$ws.Cells.Item(4, 19).Value = 1.82  # S4: 1.79 -> 1.82
$ws.Cells.Item(4, 23).Value = 1.83  # W4: 1.8 -> 1.83
$ws.Cells.Item(4, 25).Value = 18  # Y4: 1000 -> 18
$ws.Cells.Item(4, 27).Value = 80  # AA4: 1000 -> 80
$ws.Cells.Item(4, 29).Value = 9  # AC4: 1000 -> 9
$ws.Cells.Item(4, 30).Value = 16  # AD4: 1000 -> 16
$ws.Cells.Item(4, 32).Value = 21  # AF4: 1000 -> 21
$ws.Cells.Item(4, 33).Value = 12  # AG4: 1000 -> 12
$ws.Cells.Item(4, 34).Value = 16.5  # AH4: 1000 -> 16.5
$ws.Cells.Item(4, 36).Value = 48  # AJ4: 1000 -> 48
$ws.Cells.Item(4, 37).Value = 22  # AK4: 1000 -> 22
$ws.Cells.Item(4, 39).Value = 1000  # AM4: 980 -> 1000
$ws.Cells.Item(4, 40).Value = 15.5  # AN4: 1000 -> 15.5
# Row 5
$ws.Cells.Item(5, 9).Value = 4.4  # I5: 4.5 -> 4.4
$ws.Cells.Item(5, 10).Value = 3  # J5: 2.94 -> 3
$ws.Cells.Item(5, 16).Value = 1.68  # P5: 1.67 -> 1.68
$ws.Cells.Item(5, 24).Value = 980  # X5: 960 -> 980
$ws.Cells.Item(5, 25).Value = 980  # Y5: 960 -> 980
$ws.Cells.Item(5, 26).Value = 980  # Z5: 30 -> 980
$ws.Cells.Item(5, 28).Value = 10  # AB5: 960 -> 10
$ws.Cells.Item(5, 29).Value = 9.4  # AC5: 960 -> 9.4
$ws.Cells.Item(5, 30).Value = 980  # AD5: 960 -> 980
$ws.Cells.Item(5, 32).Value = 980  # AF5: 960 -> 980
$ws.Cells.Item(5, 33).Value = 980  # AG5: 960 -> 980
$ws.Cells.Item(5, 34).Value = 980  # AH5: 24 -> 980
$ws.Cells.Item(5, 35).Value = 80  # AI5: 85 -> 80
$ws.Cells.Item(5, 36).Value = 980  # AJ5: 29 -> 980
$ws.Cells.Item(5, 37).Value = 980  # AK5: 27 -> 980
$ws.Cells.Item(5, 38).Value = 60  # AL5: 65 -> 60
$ws.Cells.Item(5, 40).Value = 980  # AN5: 22 -> 980
# Row 7
$ws.Cells.Item(7, 8).Value = 2.36  # H7: 2.38 -> 2.36
$ws.Cells.Item(7, 16).Value = 1.7  # P7: 1.71 -> 1.7
$ws.Cells.Item(7, 17).Value = 2.12  # Q7: 2.02 -> 2.12
$ws.Cells.Item(7, 37).Value = 980  # AK7: 55 -> 980
# Row 8
$ws.Cells.Item(8, 12).Value = 1.32  # L8: 1.01 -> 1.32
$ws.Cells.Item(8, 17).Value = 1.89  # Q8: 1.9 -> 1.89
$ws.Cells.Item(8, 18).Value = 1.31  # R8: 1.32 -> 1.31
$ws.Cells.Item(8, 19).Value = 3.45  # S8: 3.4 -> 3.45
$ws.Cells.Item(8, 31).Value = 980  # AE8: 55 -> 980
# Row 9
$ws.Cells.Item(9, 12).Value = 1.25  # L9: 1.23 -> 1.25
$ws.Cells.Item(9, 18).Value = 1.39  # R9: 1.37 -> 1.39
$ws.Cells.Item(9, 19).Value = 2.5  # S9: 2.52 -> 2.5
$ws.Cells.Item(9, 20).Value = 1.04  # T9: 1.68 -> 1.04
# Row 10
$ws.Cells.Item(10, 8).Value = 2.76  # H10: 2.84 -> 2.76
$ws.Cells.Item(10, 10).Value = 3.8  # J10: 3.9 -> 3.8
$ws.Cells.Item(10, 14).Value = 3.8  # N10: 3.6 -> 3.8
# Row 11
$ws.Cells.Item(11, 6).Value = 5.2  # F11: 5.3 -> 5.2
$ws.Cells.Item(11, 8).Value = 1.71  # H11: 1.69 -> 1.71
$ws.Cells.Item(11, 9).Value = 1.72  # I11: 1.71 -> 1.72
$ws.Cells.Item(11, 16).Value = 2.32  # P11: 2.34 -> 2.32
$ws.Cells.Item(11, 18).Value = 1.52  # R11: 1.53 -> 1.52
$ws.Cells.Item(11, 22).Value = 2.38  # V11: 2.42 -> 2.38
$ws.Cells.Item(11, 23).Value = 1.23  # W11: 1.22 -> 1.23
$ws.Cells.Item(11, 25).Value = 11  # Y11: 10.5 -> 11
$ws.Cells.Item(11, 27).Value = 17.5  # AA11: 17 -> 17.5
$ws.Cells.Item(11, 28).Value = 21  # AB11: 23 -> 21
$ws.Cells.Item(11, 30).Value = 9.800000000000001  # AD11: 9.6 -> 9.800000000000001
$ws.Cells.Item(11, 32).Value = 42  # AF11: 44 -> 42
$ws.Cells.Item(11, 33).Value = 19.5  # AG11: 20 -> 19.5
$ws.Cells.Item(11, 36).Value = 120  # AJ11: 130 -> 120
$ws.Cells.Item(11, 37).Value = 60  # AK11: 65 -> 60
$ws.Cells.Item(11, 41).Value = 8.199999999999999  # AO11: 8 -> 8.199999999999999
# Row 12
$ws.Cells.Item(12, 13).Value = 1.07  # M12: 1.08 -> 1.07
$ws.Cells.Item(12, 16).Value = 1.96  # P12: 1.94 -> 1.96
$ws.Cells.Item(12, 36).Value = 27  # AJ12: 28 -> 27
# Row 13
$ws.Cells.Item(13, 10).Value = 4.4  # J13: 4.5 -> 4.4
$ws.Cells.Item(13, 21).Value = 2.92  # U13: 2.84 -> 2.92
$ws.Cells.Item(13, 26).Value = 40  # Z13: 38 -> 40
$ws.Cells.Item(13, 27).Value = 80  # AA13: 85 -> 80
$ws.Cells.Item(13, 31).Value = 38  # AE13: 40 -> 38
$ws.Cells.Item(13, 34).Value = 14.5  # AH13: 15 -> 14.5
$ws.Cells.Item(13, 36).Value = 22  # AJ13: 23 -> 22
$ws.Cells.Item(13, 39).Value = 48  # AM13: 50 -> 48
# Row 14
$ws.Cells.Item(14, 6).Value = 1.75  # F14: 1.73 -> 1.75
$ws.Cells.Item(14, 7).Value = 1.76  # G14: 1.75 -> 1.76
$ws.Cells.Item(14, 8).Value = 5.1  # H14: 5 -> 5.1
$ws.Cells.Item(14, 10).Value = 4.2  # J14: 4.3 -> 4.2
$ws.Cells.Item(14, 11).Value = 4.3  # K14: 4.4 -> 4.3
$ws.Cells.Item(14, 16).Value = 2.24  # P14: 2.22 -> 2.24
$ws.Cells.Item(14, 17).Value = 1.79  # Q14: 1.78 -> 1.79
$ws.Cells.Item(14, 18).Value = 1.48  # R14: 1.47 -> 1.48
$ws.Cells.Item(14, 19).Value = 2.94  # S14: 2.98 -> 2.94
$ws.Cells.Item(14, 23).Value = 2.3  # W14: 2.32 -> 2.3
$ws.Cells.Item(14, 25).Value = 19.5  # Y14: 20 -> 19.5
$ws.Cells.Item(14, 27).Value = 120  # AA14: 130 -> 120
$ws.Cells.Item(14, 30).Value = 19  # AD14: 19.5 -> 19
$ws.Cells.Item(14, 34).Value = 18  # AH14: 18.5 -> 18
$ws.Cells.Item(14, 35).Value = 60  # AI14: 65 -> 60
$ws.Cells.Item(14, 41).Value = 55  # AO14: 60 -> 55
# Row 15
$ws.Cells.Item(15, 6).Value = 1.86  # F15: 1.95 -> 1.86
$ws.Cells.Item(15, 7).Value = 2.06  # G15: 2.2 -> 2.06
$ws.Cells.Item(15, 8).Value = 4.5  # H15: 4.1 -> 4.5
$ws.Cells.Item(15, 9).Value = 5.4  # I15: 5.6 -> 5.4
$ws.Cells.Item(15, 10).Value = 3.25  # J15: 2.98 -> 3.25
$ws.Cells.Item(15, 11).Value = 3.8  # K15: 3.75 -> 3.8
$ws.Cells.Item(15, 14).Value = 3  # N15: 2.96 -> 3
$ws.Cells.Item(15, 15).Value = 1.39  # O15: 1.4 -> 1.39
$ws.Cells.Item(15, 16).Value = 1.68  # P15: 1.67 -> 1.68
$ws.Cells.Item(15, 17).Value = 2.14  # Q15: 2.16 -> 2.14
$ws.Cells.Item(15, 19).Value = 4  # S15: 3.75 -> 4
$ws.Cells.Item(15, 20).Value = 1.94  # T15: 1.92 -> 1.94
$ws.Cells.Item(15, 21).Value = 1.84  # U15: 1.86 -> 1.84
$ws.Cells.Item(15, 22).Value = 1.22  # V15: 1.23 -> 1.22
$ws.Cells.Item(15, 23).Value = 1.94  # W15: 1.83 -> 1.94
$ws.Cells.Item(15, 27).Value = 160  # AA15: 140 -> 160
$ws.Cells.Item(15, 28).Value = 9  # AB15: 10 -> 9
$ws.Cells.Item(15, 29).Value = 9.6  # AC15: 9.4 -> 9.6
$ws.Cells.Item(15, 31).Value = 95  # AE15: 85 -> 95
$ws.Cells.Item(15, 34).Value = 980  # AH15: 26 -> 980
$ws.Cells.Item(15, 35).Value = 110  # AI15: 95 -> 110
$ws.Cells.Item(15, 40).Value = 980  # AN15: 22 -> 980
$ws.Cells.Item(15, 41).Value = 130  # AO15: 110 -> 130
# Row 16
$ws.Cells.Item(16, 6).Value = 1.12  # F16: 1.09 -> 1.12
$ws.Cells.Item(16, 7).Value = 570  # G16: 2.2 -> 570
$ws.Cells.Item(16, 8).Value = 1.12  # H16: 1.04 -> 1.12
$ws.Cells.Item(16, 10).Value = 1.23  # J16: 1.09 -> 1.23
$ws.Cells.Item(16, 22).Value = 1.1  # V16: 1.22 -> 1.1
$ws.Cells.Item(16, 23).Value = 1.01  # W16: 1.83 -> 1.01
# Row 17
$ws.Cells.Item(17, 7).Value = 2.24  # G17: 2.26 -> 2.24
$ws.Cells.Item(17, 10).Value = 3.65  # J17: 3.5 -> 3.65
$ws.Cells.Item(17, 11).Value = 3.9  # K17: 3.8 -> 3.9
$ws.Cells.Item(17, 16).Value = 1.97  # P17: 1.96 -> 1.97
$ws.Cells.Item(17, 17).Value = 1.88  # Q17: 1.9 -> 1.88
$ws.Cells.Item(17, 20).Value = 1.74  # T17: 1.73 -> 1.74
$ws.Cells.Item(17, 23).Value = 1.81  # W17: 1.8 -> 1.81
$ws.Cells.Item(17, 24).Value = 18.5  # X17: 15.5 -> 18.5
